$d = $word.ActiveDocument

# Disable "smart quotes" style autocorrect so literal straight quotes survive.
$word.Options.AutoFormatAsYouTypeReplaceQuotes = $false
$word.Options.AutoFormatReplaceQuotes = $false

# ------------------------------------------------------------------
# 1. Insert a new "Meta description" paragraph right after the title.
# ------------------------------------------------------------------
$title = $d.Paragraphs.Item(1)
$title.Range.InsertParagraphAfter()

$metaPara = $d.Paragraphs.Item(2)
$metaPara.Style = "Normal"

$metaXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r><w:r><w:t xml:space="preserve">: Read our review of Cosmic Heart, a space-themed online slot game with special features and high RTP. Play for free now!</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$metaPara.Range.InsertXML($metaXml) | Out-Null

# ------------------------------------------------------------------
# 2. Remove the duplicated bold title paragraph near the end of the
#    document (it now appears twice: the Heading1 at the top and this
#    plain-text copy at the bottom).
# ------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Play Cosmic Heart Free: Review of Space-Themed Slot Game", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Collapse(0)
$rng.Find.Execute("Play Cosmic Heart Free: Review of Space-Themed Slot Game", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Expand(4) | Out-Null
$rng.Delete() | Out-Null

# ------------------------------------------------------------------
# 3. Replace the text of the italic paragraph at the very end with the
#    new image-prompt text (keeps the existing italic run formatting).
#    The phrase now also appears inside the new "Meta description"
#    paragraph, so keep searching forward until the LAST match (the
#    one in the trailing italic paragraph) is found.
# ------------------------------------------------------------------
$rng2 = $d.Content
$searchText = "Read our review of Cosmic Heart, a space-themed online slot game with special features and high RTP. Play for free now!"
$lastStart = -1
$lastEnd = -1
while ($rng2.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
    $lastStart = $rng2.Start
    $lastEnd = $rng2.End
    $rng2.Collapse(0)
}
$finalRange = $d.Range($lastStart, $lastEnd)
$finalRange.Text = 'Create a cartoon style feature image for the slot game "Cosmic Heart". The image should feature a happy Maya warrior with glasses. The Maya warrior could be seen in a spaceship or on a planet, surrounded by elements of outer space such as stars, planets, or galaxies. The image should be colorful, eye-catching, and highlight the theme of space adventure and exploration.'
